# Burndown Sprint 1 - update Day 2 ("E" column) actuals for the
# Cauldron Room work session, which ripples through the burndown
# formulas in F3:H3 (and, via the chart's live reference, the
# "Actual Burndown" series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hours logged against Day 2 for each task worked on this session.
$ws.Range("E9").Value  = 4.5   # Room A - Furnace Room Assets & Props
$ws.Range("E10").Value = 1     # Room B - Statue Room Assets & Props
$ws.Range("E11").Value = 1     # Room C - Lake Room Assets & Props
$ws.Range("E12").Value = 0.5   # Room D - Wine Cellar Assets & Props
$ws.Range("E15").Value = 0.5   # Landing - Assets & Props
$ws.Range("E16").Value = 1.5   # MSQI

# Leave the cursor where the author last left it.
$ws.Range("J23").Select()
